# Update the "想去人数" (interested count) figures that changed in this data refresh.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 3, 5, 9, 10
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 98
$ws1.Range("F5").Value = 4888
$ws1.Range("F9").Value = 743
$ws1.Range("F10").Value = 230

# Sheet "全部类型" (sheet4): rows 3, 5, 9, 11
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 98
$ws4.Range("F5").Value = 4888
$ws4.Range("F9").Value = 743
$ws4.Range("F11").Value = 230
